$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sample (SRR8423815) was sequenced between SRR8423814 and SRR8423816.
# Insert a new row right after SRR8423814 (row 20) so the table keeps its
# existing ascending sort order, then fill in its values.
$ws.Rows.Item(21).Insert() | Out-Null

$ws.Range("A21").Value = "SRR8423815"
$ws.Range("B21").Value = 1403665
$ws.Range("C21").Value = "Illumina MiSeq"

# Re-apply the table's existing ascending sort (by Sample) now that the
# range has grown by one row, so the sort range/state stays in sync.
$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add($ws.Range("A2:A121")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:D121"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Move the active selection, as it was left after the edit.
$ws.Range("A8").Select() | Out-Null
